# Updated symbol list on Mon Dec 19 07:42:48 UTC 2022 with GitHub Actions
#
# The workbook stores every data cell (column D "Price" and column E
# "Volume(1h)") as TEXT, even though most of the Price values look like
# plain numbers. Excel's COM layer auto-coerces a numeric-looking string
# assigned via Range.Value into a real Number and stamps the cell with a
# "Text" number format (so it doesn't get reinterpreted later). To keep the
# values as plain text with the default/general style (matching the
# original file), each numeric-looking cell is: (1) pre-formatted as Text
# via NumberFormat = "@", (2) assigned its new literal value, (3) reset back
# to the workbook's default "Normal" style so no stray style id is left
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Addr,
        [string]$Val
    )
    $rng = $ws.Range($Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = "Normal"
}

# --- Column D (Price) updates -------------------------------------------
Set-TextValue "D2"  "246.56"
Set-TextValue "D3"  "21.79"
Set-TextValue "D4"  "5.466"
Set-TextValue "D6"  "3.376"
Set-TextValue "D7"  "0.8021"
Set-TextValue "D9"  "0.1444"
Set-TextValue "D10" "0.07267"
Set-TextValue "D11" "0.03162"
Set-TextValue "D12" "0.02940"
Set-TextValue "D13" "0.09290"
Set-TextValue "D14" "0.001645"
Set-TextValue "D15" "3.213"
Set-TextValue "D16" "0.04717"
Set-TextValue "D17" "0.0005895"
Set-TextValue "D18" "0.006338"
Set-TextValue "D19" "0.005041"
Set-TextValue "D20" "0.001048"
Set-TextValue "D21" "0.0001501"
Set-TextValue "D22" "0.0003203"
Set-TextValue "D24" "6.425"
Set-TextValue "D25" "2.127"
Set-TextValue "D26" "0.3278"
Set-TextValue "D27" "0.1299"
Set-TextValue "D40" "0.04088"
Set-TextValue "D41" "0.006924"
Set-TextValue "D42" "0.1041"
Set-TextValue "D43" "0.002972"
Set-TextValue "D44" "0.008972"
Set-TextValue "D45" "0.00005852"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "D47" "0.7860"
Set-TextValue "D48" "0.01033"
Set-TextValue "D49" "0.00002102"
Set-TextValue "D50" "0.01011"

# --- Column E (Volume(1h)) updates ---------------------------------------
$ws.Range("E15").Value = "14MCDexMCBBestin24h"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Rows 42/43: Coin (B) and Link (C) swapped ---------------------------
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
